$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# G2 keeps its text (the start-page DSL snippet); J2:J8 used to hold "Pass" literals
# which are no longer wanted - clear them out.
$ws.Range("J2:J8").ClearContents()

# Update the selection shown when the sheet is active.
$ws.Range("J2:J8").Select()
